$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("114:114").Insert()

$ws.Range("A114").Value = 4
$ws.Range("B114").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C114").Value = 'Los Lagos'
$ws.Range("D114").Value = "2021-11-29"
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = 100112043
$ws.Range("G114").Value = 'Pepino ensalada'
$ws.Range("H114").Value = 'Sin especificar'
$ws.Range("I114").Value = 'Primera'
$ws.Range("J114").Value = 150
$ws.Range("K114").Value = 10000
$ws.Range("L114").Value = 10000
$ws.Range("M114").Value = 10000
$ws.Range("N114").Value = '$/caja 60 unidades'
$ws.Range("O114").Value = 'Región de Arica y Parinacota'
$ws.Range("P114").Value = 167
$ws.Range("Q114").Value = 60
$ws.Range("R114").Value = 'Hortaliza'
